# Apply "hybrid bold + color" highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) in specific resume bullet
# paragraphs, matching the target diff exactly.

$d = $word.ActiveDocument
$bullet = [char]8226
$pm = [char]0x00B1

# Highlight color used across all metrics: RGB 2C3E50 -> Word BGR int
$hColor = 5258796

# Each entry: exact full paragraph text (bullet + space + sentence, no
# trailing paragraph mark) plus the ordered list of substrings that must
# become bold + colored runs. Matching on the full, exact paragraph text
# avoids ambiguity between similar/overlapping bullets elsewhere in the
# document (e.g. two different "Achieved 87% ..." bullets).
$targets = @(
    @{
        Text  = "$bullet Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%"
        Terms = @("23%", "64%")
    },
    @{
        Text  = "$bullet Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from " + $pm + "4.2% to " + $pm + "2.1%"
        Terms = @("87%", "71%", ($pm + "4.2%"), ($pm + "2.1%"))
    },
    @{
        Text  = "$bullet Wrote RFP and analyzed bids from 1,200 vendors for research platform development"
        Terms = @("1,200")
    },
    @{
        Text  = "$bullet Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+"
        Terms = @('$400M', '$1B')
    },
    @{
        Text  = "$bullet Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M"
        Terms = @("73.5%", '$4.7M')
    },
    @{
        Text  = "$bullet Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%"
        Terms = @("87%", "71%")
    }
)

foreach ($target in $targets) {
    foreach ($p in $d.Paragraphs) {
        $raw = $p.Range.Text
        $trimmed = $raw.TrimEnd([char]13, [char]7)
        if ($trimmed -eq $target.Text) {
            $base = $p.Range.Start
            $searchFrom = 0
            foreach ($term in $target.Terms) {
                $idx = $trimmed.IndexOf($term, $searchFrom)
                if ($idx -ge 0) {
                    $s = $base + $idx
                    $e = $s + $term.Length
                    $sub = $d.Range($s, $e)
                    $sub.Font.Bold = 1
                    $sub.Font.Color = $hColor
                    $searchFrom = $idx + $term.Length
                }
            }
            break
        }
    }
}
